# Update "想去人数" (want-to-go count) values in column F for the
# "展览" and "全部类型" worksheets, matching the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 5
    4  = 10230
    6  = 932
    7  = 1272
    8  = 6626
    10 = 432
    13 = 3150
    14 = 34
    15 = 306
    16 = 626
    18 = 273
    19 = 273
    20 = 50
    21 = 1587
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
